$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- EDIT: update existing quantity values ---
$ws.Range("E2").Value = 20
$ws.Range("E3").Value = 26

# --- DEDUCT feature placeholder (values above represent deducted quantities) ---

# --- REPORT: append new product rows ---
$newRows = @(
    @("64000e64-3d2e-4d57-a1a8-a955d7520478", "Widget A", "A Premium Widget", 10, 15),
    @("eb3e53df-7935-439d-8c7b-04d2047f68a8", "WidgetB", "A low grade Widget", 5, 100),
    @("c7b65e32-1a9f-40f1-9b04-89b69248c9a1", "Dropdown C", "A good looking dropdown", 10, 30),
    @("5bf4ec97-f9a9-46c1-8f98-b169c5c6686d", "Dropdown D", "A premium dropdown", 50, 5)
)

$startRow = 7
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]
    $ws.Cells.Item($row, 5).Value = $data[4]
}

# Update selection to match target state
$ws.Range("L15").Select() | Out-Null
